$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (surgical in-place character edits to preserve rich-text runs) ---
$a8 = $ws.Range("A8")
$a8.Characters(21, 2).Text = "47"

$c9 = $ws.Range("C9")
$c9.Characters(27, 10).Text = "11/18/2024"
$c9.Characters(48, 10).Text = "11/24/2024"

# --- Cells that flip from a numeric value to the text placeholders "0" / "***.*" ---
# These reuse the existing style used by sibling placeholder cells (e.g. C14/F14 for "0",
# E14/H14 for "***.*"), applied via PasteSpecial so the exact existing style index is reused
# (direct NumberFormat/Value assignment would fabricate a brand-new style entry).
$zeroSrc = $ws.Range("C14")   # style 13, text "0" (shared string 20)
$naSrc = $ws.Range("E14")     # style 13, text "***.*" (shared string 21)

$zeroCells = @("D15", "C20", "D27", "F29", "G29", "F30", "G30", "D31")
foreach ($addr in $zeroCells) {
    $zeroSrc.Copy()
    $ws.Range($addr).PasteSpecial(-4122)  # xlPasteFormats
    $zeroSrc.Copy()
    $ws.Range($addr).PasteSpecial(-4163)  # xlPasteValuesAndNumberFormats
}

$naCells = @("E15", "E27", "H29", "H30", "E31")
foreach ($addr in $naCells) {
    $naSrc.Copy()
    $ws.Range($addr).PasteSpecial(-4122)  # xlPasteFormats
    $naSrc.Copy()
    $ws.Range($addr).PasteSpecial(-4163)  # xlPasteValuesAndNumberFormats
}

$excel.CutCopyMode = $false

# --- Plain numeric value updates (style untouched) ---
$ws.Range("M15").Value = -25
$ws.Range("C16").Value = 8
$ws.Range("E16").Value = 100
$ws.Range("F16").Value = 24
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = 41.176470588235
$ws.Range("I16").Value = 220
$ws.Range("J16").Value = 178
$ws.Range("K16").Value = 23.595505617977
$ws.Range("L16").Value = 25.714285714285
$ws.Range("M16").Value = -3.930131004366
$ws.Range("N16").Value = -73.137973137973
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 10
$ws.Range("E17").Value = -60
$ws.Range("F17").Value = 43
$ws.Range("G17").Value = 40
$ws.Range("H17").Value = 7.5
$ws.Range("I17").Value = 444
$ws.Range("J17").Value = 393
$ws.Range("K17").Value = 12.977099236641
$ws.Range("L17").Value = 24.022346368715
$ws.Range("M17").Value = 115.533980582524
$ws.Range("N17").Value = -17.318435754189
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 300
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 10
$ws.Range("H18").Value = 50
$ws.Range("I18").Value = 148
$ws.Range("J18").Value = 119
$ws.Range("K18").Value = 24.369747899159
$ws.Range("L18").Value = -18.681318681318
$ws.Range("M18").Value = 70.114942528735
$ws.Range("N18").Value = -65.094339622641
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = -33.333333333333
$ws.Range("G19").Value = 36
$ws.Range("H19").Value = -27.777777777777
$ws.Range("I19").Value = 346
$ws.Range("J19").Value = 395
$ws.Range("K19").Value = -12.405063291139
$ws.Range("L19").Value = -21.004566210045
$ws.Range("M19").Value = 46.610169491525
$ws.Range("N19").Value = -33.969465648855
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = -37.5
$ws.Range("J20").Value = 76
$ws.Range("K20").Value = -1.315789473684
$ws.Range("L20").Value = -5.063291139240
$ws.Range("M20").Value = 114.285714285714
$ws.Range("N20").Value = -80.366492146596
$ws.Range("C21").Value = 20
$ws.Range("E21").Value = -20
$ws.Range("F21").Value = 113
$ws.Range("G21").Value = 113
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 1258
$ws.Range("J21").Value = 1180
$ws.Range("K21").Value = 6.610169491525
$ws.Range("L21").Value = 0.318979266347
$ws.Range("M21").Value = 53.414634146341
$ws.Range("N21").Value = -54.601226993865
$ws.Range("D23").Value = 8
$ws.Range("E23").Value = -12.5
$ws.Range("F23").Value = 32
$ws.Range("G23").Value = 37
$ws.Range("H23").Value = -13.513513513513
$ws.Range("I23").Value = 371
$ws.Range("J23").Value = 382
$ws.Range("K23").Value = -2.879581151832
$ws.Range("L23").Value = -0.802139037433
$ws.Range("M23").Value = 41.603053435114
$ws.Range("C24").Value = 7
$ws.Range("D24").Value = 15
$ws.Range("E24").Value = -53.333333333333
$ws.Range("F24").Value = 49
$ws.Range("G24").Value = 76
$ws.Range("H24").Value = -35.526315789473
$ws.Range("I24").Value = 668
$ws.Range("J24").Value = 841
$ws.Range("K24").Value = -20.570749108204
$ws.Range("L24").Value = -17.530864197530
$ws.Range("M24").Value = 9.868421052631
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = -75
$ws.Range("F25").Value = 15
$ws.Range("G25").Value = 16
$ws.Range("H25").Value = -6.25
$ws.Range("I25").Value = 132
$ws.Range("J25").Value = 249
$ws.Range("K25").Value = -46.987951807228
$ws.Range("L25").Value = -30.526315789473
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 22
$ws.Range("E26").Value = -54.545454545454
$ws.Range("F26").Value = 54
$ws.Range("G26").Value = 69
$ws.Range("H26").Value = -21.739130434782
$ws.Range("I26").Value = 694
$ws.Range("J26").Value = 601
$ws.Range("K26").Value = 15.474209650582
$ws.Range("L26").Value = 40.202020202020
$ws.Range("M26").Value = 12.662337662337
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 1
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 54
$ws.Range("J28").Value = 52
$ws.Range("K28").Value = 3.846153846153
$ws.Range("L28").Value = -11.475409836065
$ws.Range("M29").Value = -62.162162162162
$ws.Range("N29").Value = -84.946236559139
$ws.Range("M30").Value = -65.625
$ws.Range("N30").Value = -87.356321839080
